$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    [double]"0.9999969141016266",
    [double]"0.9990169585194482",
    [double]"0.9999999827109688",
    [double]"0.9999999997140329",
    [double]"0.9999999958041546",
    [double]"2.88055067624043e-06",
    [double]"0.0009176260715613002",
    [double]"3.913182403822731e-09",
    [double]"2.137976670159089e-10",
    [double]"2.06349003541932e-09",
    [double]"9.999963621106398e-05",
    [double]"0.00169721851163615",
    [double]"0.9999753128130129",
    [double]"0.001769472617403162",
    [double]"67.51505815013246",
    [double]"93.11145047236465"
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
